# aggiornamento fino a 1/09/2021
# Append 9 new daily-report rows (358-366, dates 2021-08-24..2021-09-01)
# to the bottom of the existing data table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 357 (the current last row) carries the date-column style (border +
# bold + centered + custom date format) used throughout column A; copy
# just its formatting down into the new rows so the new date cells match
# the existing ones exactly, without disturbing any existing styles.
$ws.Range("A357").Copy()
$ws.Range("A358:A366").PasteSpecial(-4122)  # xlPasteFormats

$dates = 44432, 44433, 44434, 44435, 44436, 44437, 44438, 44439, 44440
$newCases = 2, 0, 0, 8, 4, 12, 4, 5, 1
$rollingSum = 18, 17, 15, 21, 21, 30, 30, 33, 34
$per100k = 70.45285529766332, 66.53880778112647, 58.71071274805276, 82.19499784727387, 82.19499784727387, 117.4214254961055, 117.4214254961055, 129.1635680457161, 133.0776155622529

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = 358 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $newCases[$i]
    $ws.Cells.Item($r, 3).Value = $rollingSum[$i]
    $ws.Cells.Item($r, 4).Value = $per100k[$i]
}
